# Update crypto price/volume table with refreshed values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.049.40"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.674.48"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").Value = "1.914.42"
$ws.Range("E12").Value = "  +3.70%  "
$ws.Range("D13").Value = "1.681.62"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "27.095.03"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").Value = "1.475.87"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.575"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.900"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.30%  "
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +10.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("D45").Value = "1.823.75"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +4.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0507"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.07%  "
